$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1489.5
$ws.Range("I19").Value = 3532.3333
$ws.Range("J19").Value = 932.36365
$ws.Range("K19").Value = 3532.3333
$ws.Range("L19").Value = 932.36365
$ws.Range("M19").Value = -3357.3333
$ws.Range("N19").Value = -1282.36365
$ws.Range("H41").Value = 642.2
$ws.Range("I41").Value = 352.75
$ws.Range("J41").Value = 1800
$ws.Range("K41").Value = 352.75
$ws.Range("L41").Value = 1800
$ws.Range("M41").Value = 87.25
$ws.Range("N41").Value = -2680
$ws.Range("H43").Value = 4377.8
$ws.Range("I43").Value = 4950
$ws.Range("J43").Value = 3996.3333
$ws.Range("K43").Value = 4950
$ws.Range("L43").Value = 3996.3333
$ws.Range("M43").Value = -4881
$ws.Range("N43").Value = -4134.3333
$ws.Range("H53").Value = 475.4
$ws.Range("I53").Value = 651.8
$ws.Range("J53").Value = 122.6
$ws.Range("K53").Value = 651.8
$ws.Range("L53").Value = 122.6
$ws.Range("M53").Value = -14.79999999999995
$ws.Range("N53").Value = -1396.6
$ws.Range("H64").Value = 3833.4443
$ws.Range("I64").Value = 3833.4443
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3833.4443
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3585.4443
$ws.Range("H67").Value = 3833.4443
$ws.Range("I67").Value = 3833.4443
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3833.4443
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2975.4443
$ws.Range("H107").Value = 1226.5555
$ws.Range("I107").Value = 1204.7142
$ws.Range("J107").Value = 1303
$ws.Range("K107").Value = 1204.7142
$ws.Range("L107").Value = 1303
$ws.Range("M107").Value = 715.2858000000001
$ws.Range("N107").Value = -5143
$ws.Range("H111").Value = 641.1111
$ws.Range("I111").Value = 348.5
$ws.Range("J111").Value = 875.2
$ws.Range("K111").Value = 1045.5
$ws.Range("L111").Value = 2625.6
$ws.Range("M111").Value = 2021.5
$ws.Range("N111").Value = -8759.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3949.6924
$ws.Range("I61").Value = 2611
$ws.Range("J61").Value = 20014
$ws.Range("K61").Value = 2611
$ws.Range("L61").Value = 20014
$ws.Range("M61").Value = -2399
$ws.Range("N61").Value = -20438
$ws.Range("H74").Value = 2649.75
$ws.Range("I74").Value = 1899.5
$ws.Range("J74").Value = 3400
$ws.Range("K74").Value = 1899.5
$ws.Range("L74").Value = 3400
$ws.Range("M74").Value = -1025.5
$ws.Range("N74").Value = -5148
$ws.Range("H77").Value = 2649.75
$ws.Range("I77").Value = 1899.5
$ws.Range("J77").Value = 3400
$ws.Range("K77").Value = 9497.5
$ws.Range("L77").Value = 17000
$ws.Range("M77").Value = -5129.5
$ws.Range("N77").Value = -25736
$ws.Range("H122").Value = 580
$ws.Range("I122").Value = 580
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1740
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 710
$ws.Range("H136").Value = 3949.6924
$ws.Range("I136").Value = 2611
$ws.Range("J136").Value = 20014
$ws.Range("K136").Value = 7833
$ws.Range("L136").Value = 60042
$ws.Range("M136").Value = -5283
$ws.Range("N136").Value = -65142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 712.5
$ws.Range("I22").Value = 528.5714
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 528.5714
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -355.5714
$ws.Range("N22").Value = -2346
$ws.Range("H86").Value = 3520.5557
$ws.Range("I86").Value = 1526.5714
$ws.Range("J86").Value = 10499.5
$ws.Range("K86").Value = 1526.5714
$ws.Range("L86").Value = 10499.5
$ws.Range("M86").Value = -403.5714
$ws.Range("N86").Value = -12745.5
$ws.Range("H89").Value = 3520.5557
$ws.Range("I89").Value = 1526.5714
$ws.Range("J89").Value = 10499.5
$ws.Range("K89").Value = 7632.857
$ws.Range("L89").Value = 52497.5
$ws.Range("M89").Value = -2016.857
$ws.Range("N89").Value = -63729.5
$ws.Range("H107").Value = 666.3333
$ws.Range("I107").Value = 666.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 666.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1253.6667
$ws.Range("H134").Value = 1415.7778
$ws.Range("I134").Value = 1665.2
$ws.Range("J134").Value = 1104
$ws.Range("K134").Value = 4995.6
$ws.Range("L134").Value = 3312
$ws.Range("M134").Value = -2460.6
$ws.Range("N134").Value = -8382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1016.8
$ws.Range("I5").Value = 1599.5
$ws.Range("J5").Value = 628.3333
$ws.Range("K5").Value = 4798.5
$ws.Range("L5").Value = 1884.9999
$ws.Range("M5").Value = -4686.5
$ws.Range("N5").Value = -2108.9999
$ws.Range("H23").Value = 1001
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 1002
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3006
$ws.Range("M23").Value = -2765
$ws.Range("N23").Value = -3476
$ws.Range("H34").Value = 874.25
$ws.Range("I34").Value = 799.4
$ws.Range("J34").Value = 999
$ws.Range("K34").Value = 2398.2
$ws.Range("L34").Value = 2997
$ws.Range("M34").Value = -2314.2
$ws.Range("N34").Value = -3165
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 1920
$ws.Range("H135").Value = 1016.8
$ws.Range("I135").Value = 1599.5
$ws.Range("J135").Value = 628.3333
$ws.Range("K135").Value = 14395.5
$ws.Range("L135").Value = 5654.9997
$ws.Range("M135").Value = -11860.5
$ws.Range("N135").Value = -10724.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1353.3636
$ws.Range("I102").Value = 1439
$ws.Range("J102").Value = 1125
$ws.Range("K102").Value = 1439
$ws.Range("L102").Value = 1125
$ws.Range("M102").Value = 183
$ws.Range("N102").Value = -4369
$ws.Range("H122").Value = 30888.334
$ws.Range("I122").Value = 40640.332
$ws.Range("J122").Value = 1632.3334
$ws.Range("K122").Value = 121920.996
$ws.Range("L122").Value = 4897.0002
$ws.Range("M122").Value = -119470.996
$ws.Range("N122").Value = -9797.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5121.75
$ws.Range("I7").Value = 5121.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5121.75
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -5009.75
$ws.Range("H22").Value = 3033.3333
$ws.Range("I22").Value = 2500
$ws.Range("J22").Value = 3140
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 3140
$ws.Range("M22").Value = -2205
$ws.Range("N22").Value = -3730
$ws.Range("H27").Value = 3033.3333
$ws.Range("I27").Value = 2500
$ws.Range("J27").Value = 3140
$ws.Range("K27").Value = 2500
$ws.Range("L27").Value = 3140
$ws.Range("M27").Value = -2393
$ws.Range("N27").Value = -3354
$ws.Range("H40").Value = 1513.75
$ws.Range("I40").Value = 1513.75
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1513.75
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1377.75
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4251
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -21256
$ws.Range("H82").Value = 1078.9375
$ws.Range("I82").Value = 1199.909
$ws.Range("J82").Value = 812.8
$ws.Range("K82").Value = 1199.909
$ws.Range("L82").Value = 812.8
$ws.Range("M82").Value = -838.9090000000001
$ws.Range("N82").Value = -1534.8
$ws.Range("H85").Value = 1078.9375
$ws.Range("I85").Value = 1199.909
$ws.Range("J85").Value = 812.8
$ws.Range("K85").Value = 1199.909
$ws.Range("L85").Value = 812.8
$ws.Range("M85").Value = 48.09099999999989
$ws.Range("N85").Value = -3308.8
$ws.Range("H122").Value = 7218.625
$ws.Range("I122").Value = 5583
$ws.Range("J122").Value = 8200
$ws.Range("K122").Value = 16749
$ws.Range("L122").Value = 24600
$ws.Range("M122").Value = -14299
$ws.Range("N122").Value = -29500
$ws.Range("H126").Value = 5121.75
$ws.Range("I126").Value = 5121.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15365.25
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -12895.25
$ws.Range("H131").Value = 55999.8
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 55999.8
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 55999.8
$ws.Range("N131").Value = -66079.8
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000.334
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 30000.334
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 30000.334
$ws.Range("N64").Value = -30496.334
$ws.Range("H67").Value = 30000.334
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 30000.334
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 30000.334
$ws.Range("N67").Value = -31716.334
$ws.Range("H86").Value = 47554.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 47554.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 47554.5
$ws.Range("N86").Value = -49800.5
$ws.Range("H89").Value = 47554.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 47554.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 237772.5
$ws.Range("N89").Value = -249004.5
$ws.Range("H126").Value = 1098.3334
$ws.Range("I126").Value = 1098.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3295.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -825.0001999999999
$ws.Range("H132").Value = 1982.9
$ws.Range("I132").Value = 1982.9
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5948.700000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3418.700000000001
